$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.555373
$ws.Range("N2").Value = 3.110746
$ws.Range("O2").Value = 0.3885765569564089
$ws.Range("P2").Value = 0.3526211663203079
$ws.Range("Q2").Value = 0.04351985499766666
$ws.Range("R2").Value = 0.261119129986
$ws.Range("S2").Value = 0.3885765569564089
$ws.Range("T2").Value = 0.3526211663203079

# Row 3
$ws.Range("O3").Value = 0.0380537990759009
$ws.Range("P3").Value = 0.05179896254485632
$ws.Range("S3").Value = 0.0380537990759009
$ws.Range("T3").Value = 0.05179896254485632

# Row 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.289563
$ws.Range("N4").Value = 0.868689
$ws.Range("O4").Value = 0.07234109989177429
$ws.Range("P4").Value = 0.09847095466798705
$ws.Range("Q4").Value = 0.008102069261000001
$ws.Range("R4").Value = 0.07291862334900001
$ws.Range("S4").Value = 0.07234109989177429
$ws.Range("T4").Value = 0.09847095466798705

# Row 5
$ws.Range("M5").Value = 1.631084
$ws.Range("N5").Value = 3.262168
$ws.Range("O5").Value = 0.4074913251205256
$ws.Range("P5").Value = 0.3697857314267338
$ws.Range("Q5").Value = 0.04563827401466666
$ws.Range("R5").Value = 0.273829644088
$ws.Range("S5").Value = 0.4074913251205256
$ws.Range("T5").Value = 0.3697857314267338

# Row 6
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04059466666666667
$ws.Range("N6").Value = 0.121784
$ws.Range("O6").Value = 0.01014170607573002
$ws.Range("P6").Value = 0.01380492528774525
$ws.Range("Q6").Value = 0.001135852304888889
$ws.Range("R6").Value = 0.010222670744
$ws.Range("S6").Value = 0.01014170607573002
$ws.Range("T6").Value = 0.01380492528774525

# Row 7
$ws.Range("M7").Value = 0.333811
$ws.Range("N7").Value = 1.001433
$ws.Range("O7").Value = 0.08339551287966027
$ws.Range("P7").Value = 0.1135182597523697
$ws.Range("Q7").Value = 0.009340143050333333
$ws.Range("R7").Value = 0.084061287453
$ws.Range("S7").Value = 0.08339551287966027
$ws.Range("T7").Value = 0.1135182597523697
